$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.556565037682828
$ws.Range("C2").Value = 0.612606004275778
$ws.Range("L2").Value = 0.599020579536134

$ws.Range("B3").Value = 0.503247309719211
$ws.Range("L3").Value = 0.605713809935328

$ws.Range("B4").Value = 0.704960018034767
$ws.Range("C4").Value = 0.780310385878786
$ws.Range("D4").Value = 0.767261516731427
$ws.Range("E4").Value = 0.757541746745335
$ws.Range("F4").Value = 0.874984008964506
$ws.Range("G4").Value = 0.847206298409435
$ws.Range("H4").Value = 0.842339547896952
$ws.Range("I4").Value = 0.65799487387268
$ws.Range("J4").Value = 0.656332785129309
$ws.Range("K4").Value = 0.703874631903231
$ws.Range("L4").Value = 0.777883926828007
$ws.Range("M4").Value = 0.92816201896394
$ws.Range("N4").Value = 0.562406199574745

$ws.Range("D5").Value = 0.330319367513297
$ws.Range("E5").Value = 0.336287925808747
$ws.Range("F5").Value = 0.557635114819894
$ws.Range("H5").Value = 0.403966853457814
$ws.Range("I5").Value = 0.352147255041012
$ws.Range("N5").Value = 0.306378388289677

$ws.Range("B6").Value = 0.410629863862209
$ws.Range("D6").Value = 0.323193155043108
$ws.Range("E6").Value = 0.444463951462934
$ws.Range("F6").Value = 0.518235877111197
$ws.Range("G6").Value = 0.41082014424506
$ws.Range("H6").Value = 0.500549487362927
$ws.Range("I6").Value = 0.453356423772496
$ws.Range("J6").Value = 0.387974052023334
$ws.Range("K6").Value = 0.32098402394466
$ws.Range("L6").Value = 0.389828757171604
$ws.Range("M6").Value = 0.665815836537034
$ws.Range("N6").Value = 0.400110573948723
